$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 272323.28
$ws.Range("I40").Value = 2030.8948
$ws.Range("J40").Value = 557631.9
$ws.Range("K40").Value = 2030.8948
$ws.Range("L40").Value = 557631.9
$ws.Range("M40").Value = -1855.8948
$ws.Range("N40").Value = -557981.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2328.6428
$ws.Range("I137").Value = 2700.1428
$ws.Range("K137").Value = 8100.428400000001
$ws.Range("M137").Value = -5550.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12824064
$ws.Range("I32").Value = 3516.5303
$ws.Range("J32").Value = 83337080
$ws.Range("K32").Value = 3516.5303
$ws.Range("L32").Value = 83337080
$ws.Range("M32").Value = -3229.5303
$ws.Range("N32").Value = -83337654

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2558
$ws.Range("I63").Value = 2349.75
$ws.Range("J63").Value = 2630.4348
$ws.Range("K63").Value = 2349.75
$ws.Range("L63").Value = 2630.4348
$ws.Range("M63").Value = -1663.75
$ws.Range("N63").Value = -4002.4348

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2558
$ws.Range("I66").Value = 2349.75
$ws.Range("J66").Value = 2630.4348
$ws.Range("K66").Value = 11748.75
$ws.Range("L66").Value = 13152.174
$ws.Range("M66").Value = -8316.75
$ws.Range("N66").Value = -20016.174

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 877.0204
$ws.Range("I74").Value = 783.6896400000001
$ws.Range("J74").Value = 1012.35
$ws.Range("K74").Value = 783.6896400000001
$ws.Range("L74").Value = 1012.35
$ws.Range("M74").Value = 90.31035999999995
$ws.Range("N74").Value = -2760.35

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 877.0204
$ws.Range("I77").Value = 783.6896400000001
$ws.Range("J77").Value = 1012.35
$ws.Range("K77").Value = 3918.4482
$ws.Range("L77").Value = 5061.75
$ws.Range("M77").Value = 449.5517999999997
$ws.Range("N77").Value = -13797.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 981250.2
$ws.Range("I132").Value = 706.2341
$ws.Range("J132").Value = 4526293.5
$ws.Range("K132").Value = 2118.7023
$ws.Range("L132").Value = 13578880.5
$ws.Range("M132").Value = 411.2977000000001
$ws.Range("N132").Value = -13583940.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2221.5356
$ws.Range("I20").Value = 2394.7
$ws.Range("J20").Value = 1788.625
$ws.Range("K20").Value = 2394.7
$ws.Range("L20").Value = 1788.625
$ws.Range("M20").Value = -2147.7
$ws.Range("N20").Value = -2282.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 20448
$ws.Range("J40").Value = 20448
$ws.Range("L40").Value = 20448
$ws.Range("N40").Value = -20978

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2452.1
$ws.Range("I94").Value = 2538
$ws.Range("K94").Value = 2538
$ws.Range("M94").Value = -2087

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3708248.8
$ws.Range("I134").Value = 1167
$ws.Range("J134").Value = 15888661
$ws.Range("K134").Value = 3501
$ws.Range("L134").Value = 47665983
$ws.Range("M134").Value = -966
$ws.Range("N134").Value = -47671053

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 27778444
$ws.Range("I58").Value = 41667320
$ws.Range("J58").Value = 696.5833
$ws.Range("K58").Value = 41667320
$ws.Range("L58").Value = 696.5833
$ws.Range("M58").Value = -41667117
$ws.Range("N58").Value = -1102.5833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 27778444
$ws.Range("I136").Value = 41667320
$ws.Range("J136").Value = 696.5833
$ws.Range("K136").Value = 125001960
$ws.Range("L136").Value = 2089.7499
$ws.Range("M136").Value = -124999410
$ws.Range("N136").Value = -7189.7499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 10000
$ws.Range("J105").Value = 10000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -35242

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 1567
$ws.Range("I110").Value = 1567
$ws.Range("K110").Value = 4701
$ws.Range("M110").Value = -611

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1412.421
$ws.Range("I129").Value = 908.75
$ws.Range("J129").Value = 1778.7273
$ws.Range("K129").Value = 2726.25
$ws.Range("L129").Value = 5336.1819
$ws.Range("M129").Value = 2273.75
$ws.Range("N129").Value = -15336.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 776.45
$ws.Range("J131").Value = 797.68475
$ws.Range("L131").Value = 2393.05425
$ws.Range("N131").Value = -12473.05425

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10618.059
$ws.Range("I70").Value = 10969.1875
$ws.Range("K70").Value = 10969.1875
$ws.Range("M70").Value = -10699.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 10618.059
$ws.Range("I73").Value = 10969.1875
$ws.Range("K73").Value = 10969.1875
$ws.Range("M73").Value = -10033.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8906.4375
$ws.Range("I132").Value = 2464
$ws.Range("J132").Value = 23079.8
$ws.Range("K132").Value = 7392
$ws.Range("L132").Value = 69239.39999999999
$ws.Range("M132").Value = -4862
$ws.Range("N132").Value = -74299.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 44780
$ws.Range("J133").Value = 44780
$ws.Range("L133").Value = 44780
$ws.Range("N133").Value = -54900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 54375
$ws.Range("J135").Value = 54375
$ws.Range("L135").Value = 54375
$ws.Range("N135").Value = -64515

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 35715692
$ws.Range("I40").Value = 1400
$ws.Range("J40").Value = 50001412
$ws.Range("K40").Value = 1400
$ws.Range("L40").Value = 50001412
$ws.Range("M40").Value = -1264
$ws.Range("N40").Value = -50001684

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1769.5714
$ws.Range("I61").Value = 1769.5714
$ws.Range("K61").Value = 1769.5714
$ws.Range("M61").Value = -1567.5714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1468.9333
$ws.Range("I68").Value = 1474.6923
$ws.Range("J68").Value = 1431.5
$ws.Range("K68").Value = 1474.6923
$ws.Range("L68").Value = 1431.5
$ws.Range("M68").Value = -725.6922999999999
$ws.Range("N68").Value = -2929.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1468.9333
$ws.Range("I71").Value = 1474.6923
$ws.Range("J71").Value = 1431.5
$ws.Range("K71").Value = 7373.461499999999
$ws.Range("L71").Value = 7157.5
$ws.Range("M71").Value = -3629.461499999999
$ws.Range("N71").Value = -14645.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1651.4166
$ws.Range("J82").Value = 4000
$ws.Range("L82").Value = 4000
$ws.Range("N82").Value = -4722

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1651.4166
$ws.Range("J85").Value = 4000
$ws.Range("L85").Value = 4000
$ws.Range("N85").Value = -6496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 25389
$ws.Range("J92").Value = 25389
$ws.Range("L92").Value = 25389
$ws.Range("N92").Value = -30381

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1769.5714
$ws.Range("I113").Value = 1769.5714
$ws.Range("K113").Value = 1769.5714
$ws.Range("M113").Value = 400.4286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 25403416
$ws.Range("I132").Value = 38097024
$ws.Range("J132").Value = 16194.533
$ws.Range("K132").Value = 114291072
$ws.Range("L132").Value = 48583.599
$ws.Range("M132").Value = -114288542
$ws.Range("N132").Value = -53643.599

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1100
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 1400
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 4200
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -8040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 125000536
$ws.Range("I113").Value = 142857660
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 428572980
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = -428570810
$ws.Range("N113").Value = -6140
